{"js": "const pairs = [[\"15+48=63\", \"47+25=72\"], [\"44+37=81\", \"90-29=61\"], [\"9+4=13\", \"33-4=29\"], [\"19+45=64\", \"72-69=3\"], [\"50-32=18\", \"80-32=48\"], [\"50-37=13\", \"68+25=93\"], [\"70-69=1\", \"17+28=45\"], [\"13+68=81\", \"54+7=61\"], [\"80-64=16\", \"66-27=39\"], [\"98-69=29\", \"79+7=86\"], [\"9+24=33\", \"36+37=73\"], [\"25+48=73\", \"24+38=62\"], [\"28+7=35\", \"37+26=63\"], [\"23+38=61\", \"57-8=49\"], [\"70-63=7\", \"58+23=81\"], [\"53+8=61\", \"75+19=94\"], [\"72-56=16\", \"29+19=48\"], [\"27+66=93\", \"41-8=33\"], [\"73-67=6\", \"83+8=91\"], [\"40-18=22\", \"38+39=77\"], [\"71-17=54\", \"4+57=61\"], [\"58+8=66\", \"38+25=63\"], [\"6+29=35\", \"71-66=5\"], [\"47+29=76\", \"33+8=41\"], [\"80-9=71\", \"19+26=45\"], [\"74-68=6\", \"50-35=15\"], [\"90-18=72\", \"83-78=5\"], [\"24+49=73\", \"72-37=35\"], [\"91-89=2\", \"19+73=92\"], [\"8+17=25\", \"20-14=6\"], [\"51-47=4\", \"80-24=56\"], [\"35+36=71\", \"43-25=18\"], [\"77+19=96\", \"57+26=83\"], [\"40-26=14\", \"54+9=63\"], [\"8+14=22\", \"29+55=84\"], [\"13+79=92\", \"29+28=57\"], [\"30-18=12\", \"7+34=41\"], [\"67-9=58\", \"61-8=53\"], [\"8+13=21\", \"55-29=26\"], [\"27+67=94\", \"18+45=63\"], [\"17+26=43\", \"90-82=8\"], [\"16+77=93\", \"75-26=49\"], [\"97-39=58\", \"71-3=68\"], [\"71-22=49\", \"38+8=46\"], [\"79+8=87\", \"87-49=38\"], [\"83-39=44\", \"27+54=81\"], [\"37+54=91\", \"56+28=84\"], [\"70-46=24\", \"53-16=37\"], [\"8+24=32\", \"55+39=94\"], [\"82-55=27\", \"4+47=51\"], [\"26+49=75\", \"7+69=76\"], [\"50-11=39\", \"59+28=87\"], [\"39+14=53\", \"44+28=72\"], [\"6+55=61\", \"10-5=5\"], [\"92-5=87\", \"46+49=95\"], [\"54-48=6\", \"71-64=7\"], [\"47+49=96\", \"82-43=39\"], [\"45+17=62\", \"62-58=4\"], [\"92-17=75\", \"22-17=5\"], [\"28+18=46\", \"15-9=6\"], [\"48+6=54\", \"86+5=91\"], [\"78+16=94\", \"47+38=85\"], [\"66+5=71\", \"84-8=76\"], [\"55-38=17\", \"3+68=71\"], [\"38+34=72\", \"10-5=5\"], [\"82-69=13\", \"72-18=54\"], [\"27+36=63\", \"80-77=3\"], [\"46-29=17\", \"96-8=88\"], [\"19+68=87\", \"18+47=65\"], [\"7+24=31\", \"83-4=79\"], [\"61-12=49\", \"8+75=83\"], [\"49+2=51\", \"47+34=81\"], [\"25+39=64\", \"85-49=36\"], [\"50-33=17\", \"34-6=28\"], [\"84+8=92\", \"74-19=55\"], [\"79+18=97\", \"5+57=62\"], [\"88-29=59\", \"78+9=87\"], [\"64-35=29\", \"51-18=33\"], [\"61-17=44\", \"41-3=38\"], [\"81-19=62\", \"92-29=63\"], [\"46+25=71\", \"25-7=18\"], [\"16+49=65\", \"82-35=47\"], [\"68+23=91\", \"60-25=35\"], [\"82-9=73\", \"4+47=51\"], [\"27+29=56\", \"60-32=28\"], [\"27-18=9\", \"30-3=27\"], [\"87+8=95\", \"37+36=73\"], [\"19+69=88\", \"8+29=37\"], [\"43+9=52\", \"54-7=47\"], [\"9+29=38\", \"78+3=81\"], [\"30-19=11\", \"29+48=77\"], [\"7+64=71\", \"49+8=57\"], [\"71-29=42\", \"51-34=17\"], [\"26-7=19\", \"18+28=46\"], [\"98-49=49\", \"38+43=81\"], [\"49+39=88\", \"34-29=5\"], [\"88+7=95\", \"14+19=33\"], [\"7+29=36\", \"36+19=55\"], [\"29+62=91\", \"88+8=96\"], [\"23-17=6\", \"91-52=39\"]];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@('15+48=63', '47+25=72')\n    ,@('44+37=81', '90-29=61')\n    ,@('9+4=13', '33-4=29')\n    ,@('19+45=64', '72-69=3')\n    ,@('50-32=18', '80-32=48')\n    ,@('50-37=13', '68+25=93')\n    ,@('70-69=1', '17+28=45')\n    ,@('13+68=81', '54+7=61')\n    ,@('80-64=16', '66-27=39')\n    ,@('98-69=29', '79+7=86')\n    ,@('9+24=33', '36+37=73')\n    ,@('25+48=73', '24+38=62')\n    ,@('28+7=35', '37+26=63')\n    ,@('23+38=61', '57-8=49')\n    ,@('70-63=7', '58+23=81')\n    ,@('53+8=61', '75+19=94')\n    ,@('72-56=16', '29+19=48')\n    ,@('27+66=93', '41-8=33')\n    ,@('73-67=6', '83+8=91')\n    ,@('40-18=22', '38+39=77')\n    ,@('71-17=54', '4+57=61')\n    ,@('58+8=66', '38+25=63')\n    ,@('6+29=35', '71-66=5')\n    ,@('47+29=76', '33+8=41')\n    ,@('80-9=71', '19+26=45')\n    ,@('74-68=6', '50-35=15')\n    ,@('90-18=72', '83-78=5')\n    ,@('24+49=73', '72-37=35')\n    ,@('91-89=2', '19+73=92')\n    ,@('8+17=25', '20-14=6')\n    ,@('51-47=4', '80-24=56')\n    ,@('35+36=71', '43-25=18')\n    ,@('77+19=96', '57+26=83')\n    ,@('40-26=14', '54+9=63')\n    ,@('8+14=22', '29+55=84')\n    ,@('13+79=92', '29+28=57')\n    ,@('30-18=12', '7+34=41')\n    ,@('67-9=58', '61-8=53')\n    ,@('8+13=21', '55-29=26')\n    ,@('27+67=94', '18+45=63')\n    ,@('17+26=43', '90-82=8')\n    ,@('16+77=93', '75-26=49')\n    ,@('97-39=58', '71-3=68')\n    ,@('71-22=49', '38+8=46')\n    ,@('79+8=87', '87-49=38')\n    ,@('83-39=44', '27+54=81')\n    ,@('37+54=91', '56+28=84')\n    ,@('70-46=24', '53-16=37')\n    ,@('8+24=32', '55+39=94')\n    ,@('82-55=27', '4+47=51')\n    ,@('26+49=75', '7+69=76')\n    ,@('50-11=39', '59+28=87')\n    ,@('39+14=53', '44+28=72')\n    ,@('6+55=61', '10-5=5')\n    ,@('92-5=87', '46+49=95')\n    ,@('54-48=6', '71-64=7')\n    ,@('47+49=96', '82-43=39')\n    ,@('45+17=62', '62-58=4')\n    ,@('92-17=75', '22-17=5')\n    ,@('28+18=46', '15-9=6')\n    ,@('48+6=54', '86+5=91')\n    ,@('78+16=94', '47+38=85')\n    ,@('66+5=71', '84-8=76')\n    ,@('55-38=17', '3+68=71')\n    ,@('38+34=72', '10-5=5')\n    ,@('82-69=13', '72-18=54')\n    ,@('27+36=63', '80-77=3')\n    ,@('46-29=17', '96-8=88')\n    ,@('19+68=87', '18+47=65')\n    ,@('7+24=31', '83-4=79')\n    ,@('61-12=49', '8+75=83')\n    ,@('49+2=51', '47+34=81')\n    ,@('25+39=64', '85-49=36')\n    ,@('50-33=17', '34-6=28')\n    ,@('84+8=92', '74-19=55')\n    ,@('79+18=97', '5+57=62')\n    ,@('88-29=59', '78+9=87')\n    ,@('64-35=29', '51-18=33')\n    ,@('61-17=44', '41-3=38')\n    ,@('81-19=62', '92-29=63')\n    ,@('46+25=71', '25-7=18')\n    ,@('16+49=65', '82-35=47')\n    ,@('68+23=91', '60-25=35')\n    ,@('82-9=73', '4+47=51')\n    ,@('27+29=56', '60-32=28')\n    ,@('27-18=9', '30-3=27')\n    ,@('87+8=95', '37+36=73')\n    ,@('19+69=88', '8+29=37')\n    ,@('43+9=52', '54-7=47')\n    ,@('9+29=38', '78+3=81')\n    ,@('30-19=11', '29+48=77')\n    ,@('7+64=71', '49+8=57')\n    ,@('71-29=42', '51-34=17')\n    ,@('26-7=19', '18+28=46')\n    ,@('98-49=49', '38+43=81')\n    ,@('49+39=88', '34-29=5')\n    ,@('88+7=95', '14+19=33')\n    ,@('7+29=36', '36+19=55')\n    ,@('29+62=91', '88+8=96')\n    ,@('23-17=6', '91-52=39')\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute(\n        $oldText,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $newText,\n        2\n    )\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
